$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.263398051261902
$ws.Range("B1").Value = 2.262205839157104
$ws.Range("C1").Value = 4.587213039398193
$ws.Range("D1").Value = 2.910417556762695
$ws.Range("E1").Value = 1.36284339427948
